$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Columns B:G contain exclusively text-like data in this sheet (inline strings).
# Excel COM auto-converts numeric-looking / percentage-looking strings to numbers,
# which would corrupt formatting (trailing zeros, thousands separators, etc). Force
# the whole data range to Text format first, write the literal strings, then restore
# the original (default/"Normal") style so the persisted cell styling is unaffected.
$dataRange = $ws.Range("B2:G51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '257.72'
$ws.Range("D3").Value = '27.66'
$ws.Range("E3").Value = '-2.74%'
$ws.Range("D4").Value = '5.233'
$ws.Range("E4").Value = '-0.06%'
$ws.Range("D5").Value = '0.05920'
$ws.Range("E5").Value = '3.88%'
$ws.Range("D6").Value = '6.692'
$ws.Range("E6").Value = '1.22%'
$ws.Range("D7").Value = '0.8698'
$ws.Range("E7").Value = '2.30%'
$ws.Range("D8").Value = '1.056'
$ws.Range("E8").Value = '21.70%'
$ws.Range("B9").Value = 'One'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D9").Value = '0.01060'
$ws.Range("E9").Value = '1,676.18%'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = '0.1415'
$ws.Range("E10").Value = '3.42%'
$ws.Range("D11").Value = '0.07197'
$ws.Range("E11").Value = '2.40%'
$ws.Range("D12").Value = '0.03261'
$ws.Range("E12").Value = '3.91%'
$ws.Range("D13").Value = '0.09220'
$ws.Range("E13").Value = '0.16%'
$ws.Range("D14").Value = '0.001551'
$ws.Range("E14").Value = '1.41%'
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").Value = '0.005954'
$ws.Range("E15").Value = '1.11%'
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").Value = '3.482'
$ws.Range("E16").Value = '-0.26%'
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D17").Value = '3.262'
$ws.Range("E17").Value = '2.06%'
$ws.Range("B18").Value = 'BTSEToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D18").Value = '2.210'
$ws.Range("E18").Value = '1.64%'
$ws.Range("B19").Value = 'BitpandaEcosystemToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D19").Value = '0.3150'
$ws.Range("E19").Value = '-0.60%'
$ws.Range("B20").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C20").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D20").Value = '0.03632'
$ws.Range("E20").Value = '10.48%'
$ws.Range("E21").Value = '-0.49%'
$ws.Range("D22").Value = '3.561'
$ws.Range("E22").Value = '1.29%'
$ws.Range("D23").Value = '0.04180'
$ws.Range("E23").Value = '2.18%'
$ws.Range("D24").Value = '0.1401'
$ws.Range("E24").Value = '1.61%'
$ws.Range("E25").Value = '0.14%'
$ws.Range("D26").Value = '0.004538'
$ws.Range("E26").Value = '9.48%'
$ws.Range("E27").Value = '0.16%'
$ws.Range("D28").Value = '0.0001940'
$ws.Range("E28").Value = '33.96%'
$ws.Range("D40").Value = '0.03824'
$ws.Range("E40").Value = '1.49%'
$ws.Range("D41").Value = '0.005496'
$ws.Range("E41").Value = '6.36%'
$ws.Range("D42").Value = '0.1106'
$ws.Range("E42").Value = '3.88%'
$ws.Range("D43").Value = '0.002302'
$ws.Range("E43").Value = '4.71%'
$ws.Range("D44").Value = '0.009938'
$ws.Range("E44").Value = '8.16%'
$ws.Range("D45").Value = '0.00005435'
$ws.Range("E45").Value = '3.11%'
$ws.Range("E46").Value = '0.15%'
$ws.Range("D47").Value = '0.1092'
$ws.Range("E47").Value = '4.06%'
$ws.Range("E49").Value = '0.15%'
$ws.Range("E50").Value = '0.15%'

# Restore original styling (removes the temporary text-format override).
$dataRange.Style = "Normal"
